$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cols = 5
$pairs = @(
    @("61×37=2257", "90×16=1440"),
    @("50×45=2250", "53×70=3710"),
    @("98×14=1372", "61×83=5063"),
    @("79×74=5846", "52×61=3172"),
    @("10×98=980", "13×53=689"),
    @("37×62=2294", "12×95=1140"),
    @("11×46=506", "20×40=800"),
    @("26×52=1352", "98×78=7644"),
    @("67×49=3283", "94×35=3290"),
    @("68×37=2516", "79×63=4977"),
    @("44×75=3300", "30×24=720"),
    @("83×79=6557", "78×43=3354"),
    @("90×19=1710", "52×99=5148"),
    @("83×89=7387", "44×48=2112"),
    @("33×21=693", "47×32=1504"),
    @("79×86=6794", "84×11=924"),
    @("92×52=4784", "10×65=650"),
    @("73×85=6205", "52×33=1716"),
    @("92×30=2760", "43×46=1978"),
    @("81×72=5832", "66×98=6468"),
    @("29×19=551", "77×62=4774"),
    @("15×11=165", "81×74=5994"),
    @("76×66=5016", "79×72=5688"),
    @("89×23=2047", "90×22=1980"),
    @("34×47=1598", "49×38=1862"),
    @("32×67=2144", "70×19=1330"),
    @("99×44=4356", "16×13=208"),
    @("77×55=4235", "53×75=3975"),
    @("54×76=4104", "25×58=1450"),
    @("64×94=6016", "53×77=4081"),
    @("50×87=4350", "62×25=1550"),
    @("92×100=9200", "48×26=1248"),
    @("83×30=2490", "19×72=1368"),
    @("62×44=2728", "87×68=5916"),
    @("70×19=1330", "91×72=6552"),
    @("79×81=6399", "17×41=697"),
    @("98×34=3332", "71×96=6816"),
    @("53×42=2226", "52×65=3380"),
    @("73×46=3358", "90×40=3600"),
    @("53×57=3021", "26×64=1664"),
    @("75×83=6225", "55×57=3135"),
    @("43×70=3010", "23×33=759"),
    @("67×93=6231", "68×70=4760"),
    @("85×83=7055", "66×54=3564"),
    @("28×22=616", "69×16=1104"),
    @("56×50=2800", "95×57=5415"),
    @("32×33=1056", "18×97=1746"),
    @("18×25=450", "98×44=4312"),
    @("25×83=2075", "22×38=836"),
    @("89×69=6141", "56×32=1792"),
    @("51×80=4080", "100×40=4000"),
    @("90×100=9000", "78×41=3198"),
    @("94×68=6392", "82×93=7626"),
    @("29×66=1914", "80×59=4720"),
    @("74×61=4514", "77×93=7161"),
    @("10×97=970", "94×21=1974"),
    @("41×83=3403", "83×67=5561"),
    @("23×70=1610", "34×42=1428"),
    @("92×98=9016", "97×39=3783"),
    @("61×74=4514", "17×92=1564"),
    @("99×83=8217", "54×69=3726"),
    @("23×90=2070", "83×27=2241"),
    @("76×53=4028", "40×59=2360"),
    @("77×22=1694", "75×15=1125"),
    @("48×17=816", "27×56=1512"),
    @("61×42=2562", "88×59=5192"),
    @("62×43=2666", "68×23=1564"),
    @("55×29=1595", "87×14=1218"),
    @("59×87=5133", "92×93=8556"),
    @("77×42=3234", "30×83=2490"),
    @("13×97=1261", "79×40=3160"),
    @("84×67=5628", "100×82=8200"),
    @("58×25=1450", "48×41=1968"),
    @("53×68=3604", "80×11=880"),
    @("69×92=6348", "33×14=462"),
    @("61×95=5795", "56×65=3640"),
    @("51×25=1275", "78×74=5772"),
    @("81×51=4131", "21×83=1743"),
    @("93×72=6696", "87×38=3306"),
    @("48×76=3648", "75×91=6825"),
    @("10×97=970", "92×73=6716"),
    @("26×76=1976", "85×12=1020"),
    @("15×29=435", "73×52=3796"),
    @("93×92=8556", "33×70=2310"),
    @("88×50=4400", "61×62=3782"),
    @("11×12=132", "31×26=806"),
    @("44×100=4400", "12×94=1128"),
    @("26×93=2418", "11×63=693"),
    @("82×61=5002", "31×22=682"),
    @("65×79=5135", "67×92=6164"),
    @("38×10=380", "34×98=3332"),
    @("32×36=1152", "21×50=1050"),
    @("95×72=6840", "34×75=2550"),
    @("41×23=943", "61×73=4453"),
    @("55×97=5335", "54×77=4158"),
    @("44×77=3388", "38×93=3534"),
    @("96×99=9504", "77×20=1540"),
    @("73×73=5329", "93×58=5394"),
    @("36×45=1620", "23×44=1012"),
    @("61×17=1037", "20×58=1160")
)
$index = 0
$mismatches = 0
foreach ($pair in $pairs) {
    $row = [math]::Floor($index / $cols) + 1
    $col = ($index % $cols) + 1
    $cell = $t.Cell($row, $col)
    $old = $pair[0]
    $new = $pair[1]
    if ($cell.Range.Text -ne ($old + "`r`a")) {
        $mismatches = $mismatches + 1
        Write-Output "MISMATCH at row $row col $col : got [$($cell.Range.Text)] expected [$old]"
    }
    $cell.Range.Text = $new
    $index = $index + 1
}
Write-Output "Updated $index cells, $mismatches mismatches"
